$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (old row 2); this shifts all subsequent rows up by one
# and drops the old last row (53), matching the new dimension A1:E52.
$ws.Rows.Item(2).Delete()

# Clear cells that must end up empty after the shift-derived values are removed
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("E10").ClearContents()

# Overwrite cells whose forecasted/lag values differ from the simple shift
$ws.Range("C9").Value = 1.173294700162053
$ws.Range("E11").Value = 1.1370912555561
$ws.Range("C12").Value = 0.9512119708358302
$ws.Range("E12").Value = 0.9990492459760025
$ws.Range("E13").Value = 0.9849212343369107
$ws.Range("C14").Value = 1.149724574326472
$ws.Range("E14").Value = 1.044407816150583
$ws.Range("E15").Value = 1.242282657891813
$ws.Range("C16").Value = 1.525861534474027
$ws.Range("E16").Value = 1.137551461271413
$ws.Range("E17").Value = 1.358051868183585
$ws.Range("C18").Value = 1.634644186146694
$ws.Range("E18").Value = 1.2772981976928
$ws.Range("E19").Value = 1.377345568933785
$ws.Range("C20").Value = 1.518308876725216
$ws.Range("E20").Value = 1.265181861560016
$ws.Range("E21").Value = 1.476362359157601
$ws.Range("C22").Value = 1.543729645060155
$ws.Range("E22").Value = 1.300276757748131
$ws.Range("C23").Value = 1.593309007378396
$ws.Range("E23").Value = 1.33496666414632
$ws.Range("C24").Value = 1.625510966833699
$ws.Range("E24").Value = 1.389880404536159
$ws.Range("E25").Value = 1.451677407676555
$ws.Range("C26").Value = 1.40478695938655
$ws.Range("E26").Value = 1.299747890163894
$ws.Range("C27").Value = 1.565661119702044
$ws.Range("E27").Value = 1.412546132271975
$ws.Range("C28").Value = 1.278716251422285
$ws.Range("E28").Value = 0.9547907592929672
$ws.Range("E29").Value = 0.6182077276742692
$ws.Range("C30").Value = 0.5944663954777107
$ws.Range("E30").Value = 1.157571511765587
$ws.Range("C31").Value = 0.287327989413555
$ws.Range("E31").Value = 0.9047322996724727
$ws.Range("E32").Value = -16.60878031022854
$ws.Range("E33").Value = -1.655311137157178
$ws.Range("C34").Value = -2.22608658996023
$ws.Range("E34").Value = 0.4768873021284703
$ws.Range("C35").Value = -1.746350382706474
$ws.Range("E35").Value = 0.7112343933969312
$ws.Range("E36").Value = 2.318221029985912
$ws.Range("E37").Value = 4.631210905746741
$ws.Range("C38").Value = 5.403124048473162
$ws.Range("E38").Value = 1.332661664932155
$ws.Range("C39").Value = 5.778434165738466
$ws.Range("E39").Value = 1.531961367047852
$ws.Range("E40").Value = 1.555532555957284
$ws.Range("E41").Value = 1.066562775371072
$ws.Range("C42").Value = 0.2475175776772698
$ws.Range("E42").Value = 1.152334833545998
$ws.Range("C43").Value = -0.2355225117835369
$ws.Range("E43").Value = 0.8847367780353999
$ws.Range("E44").Value = 0.6685736991727698
$ws.Range("E45").Value = 0.4316736535407095
$ws.Range("C46").Value = 0.2920226091170486
$ws.Range("E46").Value = 1.085023935653551
$ws.Range("C47").Value = 0.1363842982220032
$ws.Range("E47").Value = 0.9899450936446508
$ws.Range("E48").Value = 0.866426166328349
$ws.Range("E49").Value = 0.8173856700710358
$ws.Range("C50").Value = 0.7215465982331359
$ws.Range("E50").Value = 1.026748659687282
$ws.Range("C51").Value = 0.678264046940269
$ws.Range("E51").Value = 1.007646955063968
$ws.Range("E52").Value = 0.8985052439231866

Write-Output "done"
